$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update confidential note date text (A37)
$ws.Range("A37").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

# Update performance data values (columns D and E, rows 2-34)
$ws.Range("D2").Value = 0.03837597633651088
$ws.Range("E2").Value = 0.001539645881447171
$ws.Range("D3").Value = 0.02172355813531401
$ws.Range("E3").Value = -0.004050144648023224
$ws.Range("D4").Value = 0.02009290813933977
$ws.Range("E4").Value = 0.001320601381552322
$ws.Range("D5").Value = 0.03963449565285836
$ws.Range("E5").Value = -0.0007047216349542751
$ws.Range("D6").Value = 0.03711965929468199
$ws.Range("E6").Value = -0.001167769560140197
$ws.Range("D7").Value = 0.02082357985339332
$ws.Range("E7").Value = 0.002128071193654435
$ws.Range("D8").Value = 0.03816488026680344
$ws.Range("E8").Value = 0.002399650959860278
$ws.Range("D9").Value = 0.02112572117453948
$ws.Range("E9").Value = 0.005766590389016057
$ws.Range("D10").Value = 0.02627152358623763
$ws.Range("E10").Value = -0.002044571662236816
$ws.Range("D11").Value = 0.02427201946540815
$ws.Range("E11").Value = 0.0008298755186721962
$ws.Range("D12").Value = 0.05871467979772884
$ws.Range("E12").Value = 0.0009770395701025336
$ws.Range("D13").Value = 0.02647241399354192
$ws.Range("E13").Value = 0.000365230094959923
$ws.Range("D14").Value = 0.02766056795332024
$ws.Range("E14").Value = 0.0009554140127387978
$ws.Range("D15").Value = 0.03549421955793462
$ws.Range("E15").Value = -0.0237288135593221
$ws.Range("D16").Value = 0.0187100408837117
$ws.Range("E16").Value = 0.001392369813422611
$ws.Range("D17").Value = 0.03034036004433924
$ws.Range("E17").Value = -0.0119500752412145
$ws.Range("D18").Value = 0.02393700516852129
$ws.Range("E18").Value = 0.000457770656900891
$ws.Range("D19").Value = 0.1334441924665204
$ws.Range("E19").Value = -0.0006609385327166484
$ws.Range("D20").Value = 0.009404034478397341
$ws.Range("E20").Value = -0.008276405675249787
$ws.Range("D21").Value = 0.01603669447421048
$ws.Range("E21").Value = -0.0007737216009001679
$ws.Range("D22").Value = 0.01714690940163626
$ws.Range("E22").Value = -0.01134306102886362
$ws.Range("D23").Value = 0.01668625802916524
$ws.Range("E23").Value = -0.01824561403508762
$ws.Range("D24").Value = 0.02175739796328238
$ws.Range("E24").Value = -0.002441613588110347
$ws.Range("D25").Value = 0.01216047646546533
$ws.Range("E25").Value = 0.006510802012429817
$ws.Range("D26").Value = 0.04350593311447781
$ws.Range("E26").Value = 0.0005555864214679129
$ws.Range("D27").Value = 0.02517790628872015
$ws.Range("E27").Value = 0
$ws.Range("D28").Value = 0.04800679137630894
$ws.Range("E28").Value = -0.0004766444232602307
$ws.Range("D29").Value = 0.05849611318314227
$ws.Range("E29").Value = 0.00568906838580574
$ws.Range("D30").Value = 0.01346395441039935
$ws.Range("E30").Value = 0.01125827814569536
$ws.Range("D31").Value = 0.01463706844663618
$ws.Range("E31").Value = 0.002293577981651307
$ws.Range("D32").Value = 0.04412423511407131
$ws.Range("E32").Value = 0.001544799176107059
$ws.Range("D33").Value = 0.01701842548338174
$ws.Range("E33").Value = 0.006959480358800008
$ws.Range("E34").Value = -0.0008693312148640997

$ws.Protect()
